$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add sanitized name value in C2, matching B2's value ("Not applicable")
$ws.Range("C2").Value = "Not applicable"

# Leave the final selection on C5, matching the post-edit cursor position
$ws.Range("C5").Select()
